$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking value for correct answers (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update total correct marks (B12): 54 -> 90 (18 correct * 5 marks)
$ws.Range("B12").Value = 90

# Update the "score/total" summary text (E12): "50/84" -> "90/140"
$ws.Range("E12").Value = "90/140"
